$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1786.0667
$ws.Cells.Item(112, 10).Value = 1983.9231
$ws.Cells.Item(112, 12).Value = 5951.7693
$ws.Cells.Item(112, 14).Value = -8167.7693
$ws.Cells.Item(125, 8).Value = 809.1818
$ws.Cells.Item(125, 9).Value = 600
$ws.Cells.Item(125, 11).Value = 5400
$ws.Cells.Item(125, 13).Value = -2940
$ws.Cells.Item(133, 8).Value = 48236.168
$ws.Cells.Item(133, 10).Value = 48236.168
$ws.Cells.Item(133, 12).Value = 48236.168
$ws.Cells.Item(133, 14).Value = -58356.168
$ws.Cells.Item(138, 8).Value = 2133.17
$ws.Cells.Item(138, 9).Value = 1286.6111
$ws.Cells.Item(138, 10).Value = 2319
$ws.Cells.Item(138, 11).Value = 3859.8333
$ws.Cells.Item(138, 12).Value = 6957
$ws.Cells.Item(138, 13).Value = 1280.1667
$ws.Cells.Item(138, 14).Value = -17237
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2863.4614
$ws.Cells.Item(102, 9).Value = 2725
$ws.Cells.Item(102, 10).Value = 3175
$ws.Cells.Item(102, 11).Value = 2725
$ws.Cells.Item(102, 12).Value = 3175
$ws.Cells.Item(102, 13).Value = -1103
$ws.Cells.Item(102, 14).Value = -6419
$ws.Cells.Item(110, 8).Value = 1401
$ws.Cells.Item(110, 9).Value = 1334.4166
$ws.Cells.Item(110, 11).Value = 1334.4166
$ws.Cells.Item(110, 13).Value = 710.5834
$ws.Cells.Item(122, 8).Value = 1898.5
$ws.Cells.Item(122, 9).Value = 1939.1052
$ws.Cells.Item(122, 10).Value = 1744.2
$ws.Cells.Item(122, 11).Value = 5817.3156
$ws.Cells.Item(122, 12).Value = 5232.6
$ws.Cells.Item(122, 13).Value = -3367.3156
$ws.Cells.Item(122, 14).Value = -10132.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2103.5
$ws.Cells.Item(99, 9).Value = 2171.3635
$ws.Cells.Item(99, 10).Value = 2035.6364
$ws.Cells.Item(99, 11).Value = 2171.3635
$ws.Cells.Item(99, 12).Value = 2035.6364
$ws.Cells.Item(99, 13).Value = -673.3634999999999
$ws.Cells.Item(99, 14).Value = -5031.6364
$ws.Cells.Item(105, 8).Value = 3236.3333
$ws.Cells.Item(105, 9).Value = 2863.7646
$ws.Cells.Item(105, 10).Value = 4387.909
$ws.Cells.Item(105, 11).Value = 2863.7646
$ws.Cells.Item(105, 12).Value = 4387.909
$ws.Cells.Item(105, 13).Value = -1116.7646
$ws.Cells.Item(105, 14).Value = -7881.909
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1553.8182
$ws.Cells.Item(16, 9).Value = 1384.875
$ws.Cells.Item(16, 10).Value = 2004.3334
$ws.Cells.Item(16, 11).Value = 1384.875
$ws.Cells.Item(16, 12).Value = 2004.3334
$ws.Cells.Item(16, 13).Value = -1097.875
$ws.Cells.Item(16, 14).Value = -2578.3334
$ws.Cells.Item(31, 8).Value = 2470.762
$ws.Cells.Item(31, 9).Value = 2163.5
$ws.Cells.Item(31, 11).Value = 2163.5
$ws.Cells.Item(31, 13).Value = -1868.5
$ws.Cells.Item(34, 8).Value = 2470.762
$ws.Cells.Item(34, 9).Value = 2163.5
$ws.Cells.Item(34, 11).Value = 2163.5
$ws.Cells.Item(34, 13).Value = -1961.5
$ws.Cells.Item(58, 8).Value = 2600476
$ws.Cells.Item(58, 9).Value = 6995248.5
$ws.Cells.Item(58, 11).Value = 6995248.5
$ws.Cells.Item(58, 13).Value = -6995045.5
$ws.Cells.Item(92, 8).Value = 50000
$ws.Cells.Item(92, 10).Value = 50000
$ws.Cells.Item(92, 12).Value = 50000
$ws.Cells.Item(92, 14).Value = -54992
$ws.Cells.Item(105, 8).Value = 1551.75
$ws.Cells.Item(105, 9).Value = 827.625
$ws.Cells.Item(105, 11).Value = 827.625
$ws.Cells.Item(105, 13).Value = 919.375
$ws.Cells.Item(113, 8).Value = 1553.8182
$ws.Cells.Item(113, 9).Value = 1384.875
$ws.Cells.Item(113, 10).Value = 2004.3334
$ws.Cells.Item(113, 11).Value = 1384.875
$ws.Cells.Item(113, 12).Value = 2004.3334
$ws.Cells.Item(113, 13).Value = 785.125
$ws.Cells.Item(113, 14).Value = -6344.3334
$ws.Cells.Item(132, 8).Value = 2668.2104
$ws.Cells.Item(132, 9).Value = 2167.6875
$ws.Cells.Item(132, 11).Value = 6503.0625
$ws.Cells.Item(132, 13).Value = -3973.0625
$ws.Cells.Item(134, 8).Value = 2666.1372
$ws.Cells.Item(134, 9).Value = 1962
$ws.Cells.Item(134, 10).Value = 3594.318
$ws.Cells.Item(134, 11).Value = 5886
$ws.Cells.Item(134, 12).Value = 10782.954
$ws.Cells.Item(134, 13).Value = -3351
$ws.Cells.Item(134, 14).Value = -15852.954
$ws.Cells.Item(136, 8).Value = 2600476
$ws.Cells.Item(136, 9).Value = 6995248.5
$ws.Cells.Item(136, 11).Value = 20985745.5
$ws.Cells.Item(136, 13).Value = -20983195.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 69.545456
$ws.Cells.Item(38, 9).Value = 27.5
$ws.Cells.Item(38, 10).Value = 120
$ws.Cells.Item(38, 11).Value = 82.5
$ws.Cells.Item(38, 12).Value = 360
$ws.Cells.Item(38, 13).Value = 264.5
$ws.Cells.Item(38, 14).Value = -1054
$ws.Cells.Item(107, 8).Value = 1034.1
$ws.Cells.Item(107, 9).Value = 261.23077
$ws.Cells.Item(107, 10).Value = 1406.2222
$ws.Cells.Item(107, 11).Value = 783.69231
$ws.Cells.Item(107, 12).Value = 4218.6666
$ws.Cells.Item(107, 13).Value = 1136.30769
$ws.Cells.Item(107, 14).Value = -8058.6666
$ws.Cells.Item(131, 8).Value = 829.9
$ws.Cells.Item(131, 9).Value = 615.3570999999999
$ws.Cells.Item(131, 10).Value = 895.1957
$ws.Cells.Item(131, 11).Value = 1846.0713
$ws.Cells.Item(131, 12).Value = 2685.5871
$ws.Cells.Item(131, 13).Value = 3193.9287
$ws.Cells.Item(131, 14).Value = -12765.5871
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5476.0596
$ws.Cells.Item(70, 9).Value = 5296.909
$ws.Cells.Item(70, 11).Value = 5296.909
$ws.Cells.Item(70, 13).Value = -5026.909
$ws.Cells.Item(73, 8).Value = 5476.0596
$ws.Cells.Item(73, 9).Value = 5296.909
$ws.Cells.Item(73, 11).Value = 5296.909
$ws.Cells.Item(73, 13).Value = -4360.909
$ws.Cells.Item(102, 8).Value = 4243.269
$ws.Cells.Item(102, 9).Value = 3275
$ws.Cells.Item(102, 10).Value = 6871.4287
$ws.Cells.Item(102, 11).Value = 3275
$ws.Cells.Item(102, 12).Value = 6871.4287
$ws.Cells.Item(102, 13).Value = -1653
$ws.Cells.Item(102, 14).Value = -10115.4287
$ws.Cells.Item(122, 8).Value = 4145.231
$ws.Cells.Item(122, 9).Value = 4834.8887
$ws.Cells.Item(122, 10).Value = 2593.5
$ws.Cells.Item(122, 11).Value = 14504.6661
$ws.Cells.Item(122, 12).Value = 7780.5
$ws.Cells.Item(122, 13).Value = -12054.6661
$ws.Cells.Item(122, 14).Value = -12680.5
$ws.Cells.Item(126, 8).Value = 2776.923
$ws.Cells.Item(126, 9).Value = 1921
$ws.Cells.Item(126, 11).Value = 5763
$ws.Cells.Item(126, 13).Value = -3293
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6631.5386
$ws.Cells.Item(7, 9).Value = 4566.6665
$ws.Cells.Item(7, 10).Value = 8401.429
$ws.Cells.Item(7, 11).Value = 4566.6665
$ws.Cells.Item(7, 12).Value = 8401.429
$ws.Cells.Item(7, 13).Value = -4454.6665
$ws.Cells.Item(7, 14).Value = -8625.429
$ws.Cells.Item(40, 8).Value = 4757.143
$ws.Cells.Item(40, 9).Value = 4509.091
$ws.Cells.Item(40, 10).Value = 5666.6665
$ws.Cells.Item(40, 11).Value = 4509.091
$ws.Cells.Item(40, 12).Value = 5666.6665
$ws.Cells.Item(40, 13).Value = -4373.091
$ws.Cells.Item(40, 14).Value = -5938.6665
$ws.Cells.Item(126, 8).Value = 6631.5386
$ws.Cells.Item(126, 9).Value = 4566.6665
$ws.Cells.Item(126, 10).Value = 8401.429
$ws.Cells.Item(126, 11).Value = 13699.9995
$ws.Cells.Item(126, 12).Value = 25204.287
$ws.Cells.Item(126, 13).Value = -11229.9995
$ws.Cells.Item(126, 14).Value = -30144.287
$ws.Cells.Item(136, 8).Value = 4067.6365
$ws.Cells.Item(136, 9).Value = 2430.5151
$ws.Cells.Item(136, 10).Value = 8979
$ws.Cells.Item(136, 11).Value = 7291.5453
$ws.Cells.Item(136, 12).Value = 26937
$ws.Cells.Item(136, 13).Value = -4741.5453
$ws.Cells.Item(136, 14).Value = -32037
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1956.1111
$ws.Cells.Item(126, 9).Value = 1875
$ws.Cells.Item(126, 10).Value = 2605
$ws.Cells.Item(126, 11).Value = 5625
$ws.Cells.Item(126, 12).Value = 7815
$ws.Cells.Item(126, 13).Value = -3155
$ws.Cells.Item(126, 14).Value = -12755
$ws.Cells.Item(132, 8).Value = 1437.4706
$ws.Cells.Item(132, 9).Value = 726.06665
$ws.Cells.Item(132, 10).Value = 2453.762
$ws.Cells.Item(132, 11).Value = 2178.19995
$ws.Cells.Item(132, 12).Value = 7361.286
$ws.Cells.Item(132, 13).Value = 351.8000499999998
$ws.Cells.Item(132, 14).Value = -12421.286
$ws.Cells.Item(136, 8).Value = 2657.6775
$ws.Cells.Item(136, 9).Value = 2078.4062
$ws.Cells.Item(136, 11).Value = 6235.2186
$ws.Cells.Item(136, 13).Value = -3685.2186
